# Daily auto push update: insert the new 2026/01/20 09:44 UTC data point.
# The dataset is a rolling window: a new row of data is inserted at row 663
# (pushing all subsequent rows down by one), and the trailing window grows
# by one row (old D704 pair splits so the window now ends at row 705).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 663; existing rows 663:704 shift down to 664:705.
$ws.Rows.Item(663).Insert()

# Populate the newly inserted row with the new data point.
# Force column A to be treated as plain text so the date string isn't
# auto-converted into a date serial value (matches the source data which
# stores dates as literal text, e.g. "2026/12/29").
$ws.Cells.Item(663, 1).NumberFormat = "@"
$ws.Cells.Item(663, 1).Value = "2026/01/20"
$ws.Cells.Item(663, 2).Value = "火"
$ws.Cells.Item(663, 3).Value = 16
$ws.Cells.Item(663, 4).Value = 19
